# agregado FA Latam - faltan puntuales de FA y EE.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")

# --- Row 2: EGSUR / Los fogones tradicionales. Argentina / CAPS ESTRENO -> FALATAM / Breaking music / ESTRENO
$data.Range("A2").Value = "FALATAM"
$data.Range("B2").Value = "Breaking music"
$data.Range("C2").Value = "ESTRENO"
$data.Range("D2").Value = 44129.916666666664

# --- Row 3: EGNOR / Los fogones tradicionales. Argentina / BUMP -> FALATAM / Breaking music / BUMP
$data.Range("A3").Value = "FALATAM"
$data.Range("B3").Value = "Breaking music"
$data.Range("D3").Value = 44129.916666666664

# --- Row 4: EGSUR / Los fogones tradicionales. México / BUMP -> FALATAM / dowton abbey / REP
$data.Range("A4").Value = "FALATAM"
$data.Range("B4").Value = "dowton abbey"
$data.Range("C4").Value = "REP"

# --- Row 5: EE / lala land / BUMP -> FALATAM / padre brown / GEN (+ dstMex True->False)
$data.Range("A5").Value = "FALATAM"
$data.Range("B5").Value = "padre brown"
$data.Range("C5").Value = "GEN"
$data.Range("H5").Value = "False"

# --- Row 6: MCUSA / la la land usa / BUMP -> FALATAM / la la land usa / ESTRENO
$data.Range("A6").Value = "FALATAM"
$data.Range("C6").Value = "ESTRENO"

# --- Header B1: "showName" -> rich text "showName" + italic red note, with wrap text ---
$data.Range("B1").Value = "showName`nUsar la validacion para dejar indicaciones al usuario"
$data.Range("B1").WrapText = $true
$note = $data.Range("B1").Characters(10, 53)
$note.Font.Bold = $false
$note.Font.Italic = $true
$note.Font.Size = 12
$note.Font.Color = 255
$note.Font.Name = "Calibri (Cuerpo)"

# --- genDateStr column (E) for the FALATAM rows ---
$data.Range("E2").Value = "no"
$data.Range("E3").Value = "no"

# --- Selection moved from B13 to B6 ---
$data.Range("B6").Select()
